$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C9").Value = "Swiss TPH & U Basel"
$ws.Range("C42").Value = "DSM Firmenich"
$ws.Range("C2").Value = "Merck"
$ws.Range("C12").Value = "Cogitamen"
$ws.Range("C21").Value = "Swiss TPH & U Basel"
$ws.Range("C44").Value = "Swiss TPH & U Basel"
$ws.Range("C44").Font.Name = "Calibri"
$ws.Range("C44").Font.Size = 11
$ws.Range("C44").Font.Bold = $false
$ws.Range("C44").Font.Color = 0

$ws.Range("C13").Select()
